# Workbook was re-uploaded with fewer data rows and an updated A2 value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 value changed from 85 to 55 (B2 stays 29).
$ws.Range("A2").Value = 55

# Rows 3-8 were removed entirely, shrinking the used range to A1:B2.
$ws.Range("A3:B8").EntireRow.Delete()

# Selection moved from A3 to B7.
$ws.Range("B7").Select()
